$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 68.443746
$ws.Cells.Item(2, 8).Value = 205.331238
$ws.Cells.Item(2, 9).Value = 0.1596169534001499
$ws.Cells.Item(2, 10).Value = 0.1596169534001499
$ws.Cells.Item(2, 13).Value = 0.5550926666666666
$ws.Cells.Item(2, 14).Value = 1.665278
$ws.Cells.Item(2, 15).Value = 0.1208967663154349
$ws.Cells.Item(2, 16).Value = 0.1208967663154349
$ws.Cells.Item(2, 17).Value = 37.99262148379599
$ws.Cells.Item(2, 18).Value = 341.933593354164
$ws.Cells.Item(2, 19).Value = 0.01929717351519958
$ws.Cells.Item(2, 20).Value = 0.01929717351519958
$ws.Cells.Item(3, 7).Value = 68.443746
$ws.Cells.Item(3, 8).Value = 205.331238
$ws.Cells.Item(3, 9).Value = 0.1596169534001499
$ws.Cells.Item(3, 10).Value = 0.1596169534001499
$ws.Cells.Item(3, 15).Value = 0.7377399926530269
$ws.Cells.Item(3, 16).Value = 0.7377399926530268
$ws.Cells.Item(3, 17).Value = 231.839751786202
$ws.Cells.Item(3, 18).Value = 2086.557766075818
$ws.Cells.Item(3, 19).Value = 0.1177558100287251
$ws.Cells.Item(3, 20).Value = 0.1177558100287251
$ws.Cells.Item(4, 7).Value = 68.443746
$ws.Cells.Item(4, 8).Value = 205.331238
$ws.Cells.Item(4, 9).Value = 0.1596169534001499
$ws.Cells.Item(4, 10).Value = 0.1596169534001499
$ws.Cells.Item(4, 13).Value = 0.5311786666666667
$ws.Cells.Item(4, 14).Value = 1.593536
$ws.Cells.Item(4, 15).Value = 0.1156884012202364
$ws.Cells.Item(4, 16).Value = 0.1156884012202364
$ws.Cells.Item(4, 17).Value = 36.355857741952
$ws.Cells.Item(4, 18).Value = 327.202719677568
$ws.Cells.Item(4, 19).Value = 0.01846583014650833
$ws.Cells.Item(4, 20).Value = 0.01846583014650832
$ws.Cells.Item(5, 7).Value = 68.443746
$ws.Cells.Item(5, 8).Value = 205.331238
$ws.Cells.Item(5, 9).Value = 0.1596169534001499
$ws.Cells.Item(5, 10).Value = 0.1596169534001499
$ws.Cells.Item(5, 13).Value = 0.117885
$ws.Cells.Item(5, 14).Value = 0.353655
$ws.Cells.Item(5, 15).Value = 0.02567483981130185
$ws.Cells.Item(5, 16).Value = 0.02567483981130185
$ws.Cells.Item(5, 17).Value = 8.068490997210001
$ws.Cells.Item(5, 18).Value = 72.61641897489001
$ws.Cells.Item(5, 19).Value = 0.004098139709716882
$ws.Cells.Item(5, 20).Value = 0.004098139709716881
$ws.Cells.Item(6, 9).Value = 0.4159650732941736
$ws.Cells.Item(6, 10).Value = 0.4159650732941736
$ws.Cells.Item(6, 13).Value = 0.5550926666666666
$ws.Cells.Item(6, 14).Value = 1.665278
$ws.Cells.Item(6, 15).Value = 0.1208967663154349
$ws.Cells.Item(6, 16).Value = 0.1208967663154349
$ws.Cells.Item(6, 17).Value = 99.00955533543065
$ws.Cells.Item(6, 18).Value = 891.0859980188759
$ws.Cells.Item(6, 19).Value = 0.05028883226142845
$ws.Cells.Item(6, 20).Value = 0.05028883226142845
$ws.Cells.Item(7, 9).Value = 0.4159650732941736
$ws.Cells.Item(7, 10).Value = 0.4159650732941736
$ws.Cells.Item(7, 15).Value = 0.7377399926530269
$ws.Cells.Item(7, 16).Value = 0.7377399926530268
$ws.Cells.Item(7, 19).Value = 0.3068740701159595
$ws.Cells.Item(7, 20).Value = 0.3068740701159594
$ws.Cells.Item(8, 9).Value = 0.4159650732941736
$ws.Cells.Item(8, 10).Value = 0.4159650732941736
$ws.Cells.Item(8, 13).Value = 0.5311786666666667
$ws.Cells.Item(8, 14).Value = 1.593536
$ws.Cells.Item(8, 15).Value = 0.1156884012202364
$ws.Cells.Item(8, 16).Value = 0.1156884012202364
$ws.Cells.Item(8, 17).Value = 94.74411525943466
$ws.Cells.Item(8, 18).Value = 852.697037334912
$ws.Cells.Item(8, 19).Value = 0.04812233429286141
$ws.Cells.Item(8, 20).Value = 0.0481223342928614
$ws.Cells.Item(9, 9).Value = 0.4159650732941736
$ws.Cells.Item(9, 10).Value = 0.4159650732941736
$ws.Cells.Item(9, 13).Value = 0.117885
$ws.Cells.Item(9, 14).Value = 0.353655
$ws.Cells.Item(9, 15).Value = 0.02567483981130185
$ws.Cells.Item(9, 16).Value = 0.02567483981130185
$ws.Cells.Item(9, 17).Value = 21.02665398339
$ws.Cells.Item(9, 18).Value = 189.23988585051
$ws.Cells.Item(9, 19).Value = 0.01067983662392434
$ws.Cells.Item(9, 20).Value = 0.01067983662392434
$ws.Cells.Item(10, 7).Value = 88.88346833333333
$ws.Cells.Item(10, 8).Value = 266.650405
$ws.Cells.Item(10, 9).Value = 0.2072842188241036
$ws.Cells.Item(10, 10).Value = 0.2072842188241036
$ws.Cells.Item(10, 13).Value = 0.5550926666666666
$ws.Cells.Item(10, 14).Value = 1.665278
$ws.Cells.Item(10, 15).Value = 0.1208967663154349
$ws.Cells.Item(10, 16).Value = 0.1208967663154349
$ws.Cells.Item(10, 17).Value = 49.33856145973221
$ws.Cells.Item(10, 18).Value = 444.0470531375899
$ws.Cells.Item(10, 19).Value = 0.02505999176405512
$ws.Cells.Item(10, 20).Value = 0.02505999176405511
$ws.Cells.Item(11, 7).Value = 88.88346833333333
$ws.Cells.Item(11, 8).Value = 266.650405
$ws.Cells.Item(11, 9).Value = 0.2072842188241036
$ws.Cells.Item(11, 10).Value = 0.2072842188241036
$ws.Cells.Item(11, 15).Value = 0.7377399926530269
$ws.Cells.Item(11, 16).Value = 0.7377399926530268
$ws.Cells.Item(11, 17).Value = 301.0752981915505
$ws.Cells.Item(11, 18).Value = 2709.677683723955
$ws.Cells.Item(11, 19).Value = 0.1529218580723826
$ws.Cells.Item(11, 20).Value = 0.1529218580723826
$ws.Cells.Item(12, 7).Value = 88.88346833333333
$ws.Cells.Item(12, 8).Value = 266.650405
$ws.Cells.Item(12, 9).Value = 0.2072842188241036
$ws.Cells.Item(12, 10).Value = 0.2072842188241036
$ws.Cells.Item(12, 13).Value = 0.5311786666666667
$ws.Cells.Item(12, 14).Value = 1.593536
$ws.Cells.Item(12, 15).Value = 0.1156884012202364
$ws.Cells.Item(12, 16).Value = 0.1156884012202364
$ws.Cells.Item(12, 17).Value = 47.21300219800889
$ws.Cells.Item(12, 18).Value = 424.91701978208
$ws.Cells.Item(12, 19).Value = 0.02398037987394618
$ws.Cells.Item(12, 20).Value = 0.02398037987394618
$ws.Cells.Item(13, 7).Value = 88.88346833333333
$ws.Cells.Item(13, 8).Value = 266.650405
$ws.Cells.Item(13, 9).Value = 0.2072842188241036
$ws.Cells.Item(13, 10).Value = 0.2072842188241036
$ws.Cells.Item(13, 13).Value = 0.117885
$ws.Cells.Item(13, 14).Value = 0.353655
$ws.Cells.Item(13, 15).Value = 0.02567483981130185
$ws.Cells.Item(13, 16).Value = 0.02567483981130185
$ws.Cells.Item(13, 17).Value = 10.478027664475
$ws.Cells.Item(13, 18).Value = 94.302248980275
$ws.Cells.Item(13, 19).Value = 0.0053219891137197
$ws.Cells.Item(13, 20).Value = 0.005321989113719699
$ws.Cells.Item(14, 7).Value = 93.106949
$ws.Cells.Item(14, 8).Value = 279.320847
$ws.Cells.Item(14, 9).Value = 0.2171337544815728
$ws.Cells.Item(14, 10).Value = 0.2171337544815728
$ws.Cells.Item(14, 13).Value = 0.5550926666666666
$ws.Cells.Item(14, 14).Value = 1.665278
$ws.Cells.Item(14, 15).Value = 0.1208967663154349
$ws.Cells.Item(14, 16).Value = 0.1208967663154349
$ws.Cells.Item(14, 17).Value = 51.68298460560732
$ws.Cells.Item(14, 18).Value = 465.146861450466
$ws.Cells.Item(14, 19).Value = 0.02625076877475172
$ws.Cells.Item(14, 20).Value = 0.02625076877475172
$ws.Cells.Item(15, 7).Value = 93.106949
$ws.Cells.Item(15, 8).Value = 279.320847
$ws.Cells.Item(15, 9).Value = 0.2171337544815728
$ws.Cells.Item(15, 10).Value = 0.2171337544815728
$ws.Cells.Item(15, 15).Value = 0.7377399926530269
$ws.Cells.Item(15, 16).Value = 0.7377399926530268
$ws.Cells.Item(15, 17).Value = 315.3815097398463
$ws.Cells.Item(15, 18).Value = 2838.433587658617
$ws.Cells.Item(15, 19).Value = 0.1601882544359597
$ws.Cells.Item(15, 20).Value = 0.1601882544359597
$ws.Cells.Item(16, 7).Value = 93.106949
$ws.Cells.Item(16, 8).Value = 279.320847
$ws.Cells.Item(16, 9).Value = 0.2171337544815728
$ws.Cells.Item(16, 10).Value = 0.2171337544815728
$ws.Cells.Item(16, 13).Value = 0.5311786666666667
$ws.Cells.Item(16, 14).Value = 1.593536
$ws.Cells.Item(16, 15).Value = 0.1156884012202364
$ws.Cells.Item(16, 16).Value = 0.1156884012202364
$ws.Cells.Item(16, 17).Value = 49.45642502722134
$ws.Cells.Item(16, 18).Value = 445.1078252449921
$ws.Cells.Item(16, 19).Value = 0.02511985690692051
$ws.Cells.Item(16, 20).Value = 0.02511985690692051
$ws.Cells.Item(17, 7).Value = 93.106949
$ws.Cells.Item(17, 8).Value = 279.320847
$ws.Cells.Item(17, 9).Value = 0.2171337544815728
$ws.Cells.Item(17, 10).Value = 0.2171337544815728
$ws.Cells.Item(17, 13).Value = 0.117885
$ws.Cells.Item(17, 14).Value = 0.353655
$ws.Cells.Item(17, 15).Value = 0.02567483981130185
$ws.Cells.Item(17, 16).Value = 0.02567483981130185
$ws.Cells.Item(17, 17).Value = 10.975912682865
$ws.Cells.Item(17, 18).Value = 98.78321414578501
$ws.Cells.Item(17, 19).Value = 0.005574874363940929
$ws.Cells.Item(17, 20).Value = 0.005574874363940928
